$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.740.11"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "3.502.24"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'593.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").Value = "'169.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.60%  "
$ws.Range("E9").Value = "  +8.64%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "4.111.27"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "'28.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("E15").Value = "  +3.88%  "
$ws.Range("D16").Value = "66.766.72"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "3.501.54"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "'14.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").Value = "'397.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.34%  "
$ws.Range("D21").Value = "'7.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "'73.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D26").Value = "'10.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").Value = "'6.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("D32").Value = "'23.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.24%  "
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").Value = "'1.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.14%  "
$ws.Range("D35").Value = "'162.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").Value = "'0.904"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").Value = "'1.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("D39").Value = "'4.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("D40").Value = "'0.0747"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").Value = "'26.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("D42").Value = "'27.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.75%  "
$ws.Range("D43").Value = "2.804.92"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'42.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("D47").Value = "'343.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("E48").Value = "  +3.24%  "
$ws.Range("D49").Value = "'34.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.22%  "
$ws.Range("D50").Value = "'0.858"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("E51").Value = "  +2.51%  "
